$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '42.997.71'
$ws.Cells.Item(2, 5).Value = '  +0.63%  '
$ws.Cells.Item(3, 4).Value = '2.302.74'
$ws.Cells.Item(3, 5).Value = '  -0.35%  '
$ws.Cells.Item(4, 5).Value = '  +0.26%  '
$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '309.88'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -2.52%  '
$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '104.60'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -0.63%  '
$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.627'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.50%  '
$ws.Cells.Item(8, 5).Value = '  +0.14%  '
$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.605'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -0.44%  '
$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '39.58'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -1.27%  '
$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0904'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -0.58%  '
$ws.Cells.Item(12, 5).Value = '  -3.36%  '
$ws.Cells.Item(13, 5).Value = '  +0.15%  '
$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.992'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +1.41%  '
$ws.Cells.Item(15, 4).Value = '2.783.21'
$ws.Cells.Item(15, 5).Value = '  +4.75%  '
$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '15.30'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -0.77%  '
$ws.Cells.Item(17, 4).Value = '2.296.04'
$ws.Cells.Item(17, 5).Value = '  -0.52%  '
$ws.Cells.Item(18, 4).Value = '42.832.36'
$ws.Cells.Item(18, 5).Value = '  +0.40%  '
$ws.Cells.Item(19, 5).Value = '  -4.21%  '
$ws.Cells.Item(20, 2).Value = 'ShibaInu'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0000105'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -0.97%  '
$ws.Cells.Item(21, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$cell = $ws.Cells.Item(21, 4)
$cell.NumberFormat = "@"
$cell.Value = '13.65'
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +0.46%  '
$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '73.35'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.93%  '
$ws.Cells.Item(23, 5).Value = '  -3.33%  '
$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '266.80'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -1.27%  '
$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.23'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -0.71%  '
$ws.Cells.Item(26, 5).Value = '  +0.29%  '
$ws.Cells.Item(27, 2).Value = 'Cosmos'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '10.95'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +0.56%  '
$ws.Cells.Item(28, 2).Value = 'Filecoin'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.34'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  +16.75%  '
$ws.Cells.Item(29, 5).Value = '  -1.11%  '
$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '22.25'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -1.78%  '
$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '36.24'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -5.71%  '
$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '164.83'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -1.34%  '
$ws.Cells.Item(33, 5).Value = '  -3.23%  '
$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.66'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  +2.57%  '
$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.131'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -0.81%  '
$ws.Cells.Item(36, 5).Value = '  -3.12%  '
$ws.Cells.Item(37, 5).Value = '  -1.17%  '
$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.0348'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -1.51%  '
$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.84'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +2.81%  '
$ws.Cells.Item(40, 5).Value = '  -2.99%  '
$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '107.79'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +7.76%  '
$ws.Cells.Item(42, 5).Value = '  -3.84%  '
$cell = $ws.Cells.Item(43, 4)
$cell.NumberFormat = "@"
$cell.Value = '71.13'
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  +1.06%  '
$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.227'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.18%  '
$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.01'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.34%  '
$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '12.16'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -2.15%  '
$ws.Cells.Item(47, 4).Value = '1.727.61'
$ws.Cells.Item(47, 5).Value = '  +5.43%  '
$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '111.04'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -4.86%  '
$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '76.71'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -6.89%  '
$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.69'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -2.07%  '
$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.14'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -2.97%  '
